# Revision de solicitud y aprobacion de peticion
# Fills in the "Informacion llenada por Direccion" / "Cierre del Cambio"
# sections plus the real completion date of the first activity row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fecha Real for the first "Actividades a Realizar" row (16-jun-2015)
$ws.Range("G13").Value = 42171

# Informacion llenada por Direccion
$ws.Range("C35").Value = "44.62"
$ws.Range("C36").Value = "No aplica"
$ws.Range("C37").Value = "Se incorpora el cambio, el costo lo absorbe la empresa y no hay cambio en fecha de entrega"
$ws.Rows("37").RowHeight = 25.5
$ws.Range("C38").Value = 42171
$ws.Range("C39").Value = 42171
$ws.Range("C40").Value = 42171

# Cierre del Cambio
$ws.Range("C42").Value = "Fidel Reyna"
$ws.Range("C43").Value = 42171
$ws.Range("C44").Value = "Aprobado"

# Reviewer scrolled down to review/approve row 36 before saving
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("C36").Select()
